# Clear the per-lap tracking data (columns C:G, rows 2:21) on the "Score"
# sheet. The header row (row 1) and the Team/Name columns (A:B) are left
# untouched; only the numeric lap/rep data and the computed G-column
# formula are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")

$ws.Range("C2:G21").Clear()

# Match the author's final selection in the saved file.
$ws.Range("J10").Select()
